$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full attribution data: the channel list is now sorted alphabetically and
# several new media channels were added. Row 1 (headers) is unchanged;
# rows 2-17 hold the (reshuffled + expanded) per-channel data.
$data = @(
  @{ Row = 2;  A = "amazon_media_cost";        B = 325.4689352158736;  C = 947.5823100035;      D = 0.3434729962557779 },
  @{ Row = 3;  A = "bingsearch_media_cost";     B = 18604.99254492222;  C = 71190.62;             D = 0.2613405044782897 },
  @{ Row = 4;  A = "ctv_media_cost";            B = 1370.421463794661;  C = 932751.4633993;       D = 0.001469224672990944 },
  @{ Row = 5;  A = "criteo_media_cost";         B = 34.99542947131734;  C = 145.4036846402;       D = 0.240677734941265 },
  @{ Row = 6;  A = "dv360_media_cost";          B = 618.4626357864217;  C = 271129.18;            D = 0.002281062612981833 },
  @{ Row = 7;  A = "facebook_media_cost";       B = 31233.80555115086;  C = 3079426.160085;       D = 0.01014273566809237 },
  @{ Row = 8;  A = "googlesearch_media_cost";   B = 32679.13737021944;  C = 1494447.761988;       D = 0.02186703222516643 },
  @{ Row = 9;  A = "influential_media_cost";    B = 1282.53220170682;   C = 109551.439979;        D = 0.01170712317385029 },
  @{ Row = 10; A = "lineartv_media_cost";       B = 1445.577637401062;  C = 1467680;              D = 0.0009849406119869875 },
  @{ Row = 11; A = "pinterest_media_cost";      B = 3661.441939006166;  C = 21588.25;             D = 0.1696034620224504 },
  @{ Row = 12; A = "radio_media_cost";          B = 0;                  C = 0;                    D = $null },
  @{ Row = 13; A = "snapchat_media_cost";       B = 727.1401467448559;  C = 957.659705;           D = 0.7592886522722138 },
  @{ Row = 14; A = "thetradedesk_media_cost";   B = 1639.986963078731;  C = 133960.5438629989;    D = 0.0122423134139851 },
  @{ Row = 15; A = "tinder_media_cost";         B = 176.9075242613488;  C = 56425.47;             D = 0.003135242369471603 },
  @{ Row = 16; A = "twitch_media_cost";         B = 50.96558290840737;  C = 12709.62;             D = 0.004010000527821239 },
  @{ Row = 17; A = "youtube_media_cost";        B = 29554.89816805543;  C = 229054.302767;        D = 0.1290300937857493 }
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    if ($null -ne $row.D) {
        $ws.Cells.Item($r, 4).Value = $row.D
    }
}

# Rows 2-6 already carry the bold/border/centered label style from the
# original file. Stamp the same formatting onto the newly-added rows
# (7-17) by copying it from an existing label cell.
$ws.Range("A2").Copy()
$ws.Range("A7:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
